$wb = $excel.ActiveWorkbook

# --- Update header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Reuse the header formatting (bold, centered, bordered) from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse the date-column number formatting from column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A81").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$forecastRows = @(
    @(44934.99999999999, 228, -1527.259782589761, 2022.401680785223),
    @(44941.99999999999, 250, -1585.203210458682, 1954.322133764208),
    @(44948.99999999999, 271, -1497.655306065355, 2085.692818848449),
    @(44955.99999999999, 293, -1583.223926276966, 1985.00095960814),
    @(44962.99999999999, 314, -1438.909896596134, 2072.442272331328),
    @(44983.99999999999, 379, -1369.261460358188, 2175.99736853452),
    @(44990.99999999999, 401, -1237.784866710784, 2203.904643081566),
    @(44997.99999999999, 422, -1245.977495907375, 2277.793815089445),
    @(45011.99999999999, 465, -1198.208086466249, 2212.66029310278),
    @(45018.99999999999, 487, -1333.350385121826, 2437.904314467292),
    @(45025.99999999999, 508, -1101.764959205173, 2286.077369983473),
    @(45032.99999999999, 530, -1180.961397094168, 2288.514873610382),
    @(45039.99999999999, 552, -1223.925826488454, 2463.045052996644),
    @(45046.99999999999, 573, -1133.965287057402, 2383.692773590765),
    @(45053.99999999999, 595, -1075.703826688562, 2408.634670524624),
    @(45060.99999999999, 616, -1104.756268494925, 2416.316974088471),
    @(45067.99999999999, 638, -1174.507857015457, 2338.32068711844),
    @(45074.99999999999, 659, -1129.716880968276, 2414.527609582547),
    @(45081.99999999999, 681, -1044.066292832892, 2420.063035655825),
    @(45088.99999999999, 703, -879.8830952969346, 2490.52183995366),
    @(45095.99999999999, 724, -1115.366766295619, 2398.151738184424),
    @(45102.99999999999, 746, -1058.585896767193, 2402.581870321799),
    @(45109.99999999999, 767, -1093.838664489378, 2466.507800991703),
    @(45116.99999999999, 789, -922.3981950985269, 2657.612695936684),
    @(45123.99999999999, 810, -1100.974478285316, 2637.796588794017),
    @(45130.99999999999, 832, -975.6147510478107, 2455.01551704653),
    @(45137.99999999999, 854, -976.5287285466106, 2582.971259029028),
    @(45144.99999999999, 875, -873.8291442066546, 2736.258341217553),
    @(45151.99999999999, 897, -982.6864648678048, 2623.196066371137),
    @(45158.99999999999, 918, -971.5986648792223, 2657.13898030842),
    @(45172.99999999999, 961, -814.6566541336845, 2627.704399610978),
    @(45179.99999999999, 983, -728.0196900899999, 2647.431197222932),
    @(45186.99999999999, 1005, -812.9953602347631, 2732.081905018562),
    @(45193.99999999999, 1026, -706.1690271562345, 2634.891882804784),
    @(45200.99999999999, 1048, -719.6585743772013, 2676.378086098766),
    @(45207.99999999999, 1069, -784.1199586061446, 2923.158034167509),
    @(45214.99999999999, 1091, -734.8634029123752, 2893.666031542527),
    @(45221.99999999999, 1112, -715.7578291285679, 3046.992367455227),
    @(45228.99999999999, 1134, -732.494890699954, 2908.716943373037),
    @(45235.99999999999, 1156, -694.5034907710918, 2841.504645748979),
    @(45242.99999999999, 1177, -611.1147492329956, 2940.464669833027),
    @(45249.99999999999, 1199, -518.4952236575684, 3003.787983901716),
    @(45263.99999999999, 1242, -507.3027945773366, 3088.444601776509),
    @(45270.99999999999, 1263, -526.7150626120736, 3080.083693120826),
    @(45277.99999999999, 1285, -447.141643431487, 3118.574527289038),
    @(45298.99999999999, 1350, -360.5170763840422, 3149.435165964163),
    @(45312.99999999999, 1393, -240.8868880377538, 3126.672066149968),
    @(45326.99999999999, 1436, -310.115792288789, 3121.063145690302),
    @(45396.99999999999, 1652, -34.93838985389919, 3402.368048248141),
    @(45410.99999999999, 1695, -11.57956975653138, 3473.89662329242),
    @(45424.99999999999, 1738, 56.74992056956052, 3392.018101072607),
    @(45431.99999999999, 1760, -54.74645631245239, 3736.060700747733),
    @(45438.99999999999, 1781, 101.932841951017, 3547.397511449508),
    @(45445.99999999999, 1803, 48.68718517799168, 3711.923219944881),
    @(45459.99999999999, 1846, 129.2866984833626, 3545.487221479386),
    @(45466.99999999999, 1867, 101.4809383241974, 3606.763648990505),
    @(45473.99999999999, 1889, 90.14131879791232, 3527.936575344053),
    @(45487.99999999999, 1932, 156.8976092064261, 3719.296188449637),
    @(45494.99999999999, 1954, 191.773833828613, 3820.525763930658),
    @(45508.99999999999, 1997, 267.9537809749264, 3754.077032854536),
    @(45515.99999999999, 2018, 204.3969606522896, 3682.853925084066),
    @(45522.99999999999, 2040, 332.5109432260878, 3801.982147971204),
    @(45529.99999999999, 2062, 283.6470070232545, 3808.419024196115),
    @(45543.99999999999, 2105, 308.4687027990515, 3891.338238619043),
    @(45550.99999999999, 2126, 480.7509874053063, 3876.722498612969),
    @(45557.99999999999, 2148, 418.2077135670474, 3925.806948507016),
    @(45564.99999999999, 2169, 402.5439100758326, 4011.788035690754),
    @(45571.99999999999, 2191, 463.8460885578023, 3977.09590536891),
    @(45585.99999999999, 2234, 489.6136679254113, 4154.467860958982),
    @(45592.99999999999, 2256, 431.6568327304701, 4059.176630309345),
    @(45599.99999999999, 2277, 555.1259433064483, 4036.622097236363),
    @(45634.99999999999, 2385, 546.8867293423654, 4075.372384505884),
    @(45641.99999999999, 2407, 537.340280882801, 4200.654035333809),
    @(45648.99999999999, 2428, 686.5052666017937, 4077.558430994822),
    @(45655.99999999999, 2450, 739.2476615599994, 4233.819636629078),
    @(45662.99999999999, 2471, 591.4493853546702, 4097.264039295113),
    @(45669.99999999999, 2493, 744.8088334443081, 4202.57701224743),
    @(45676.99999999999, 2515, 686.4949647378354, 4188.683702369289),
    @(45683.99999999999, 2536, 795.5943321644613, 4306.034302376303),
    @(45690.99999999999, 2558, 821.4519194078587, 4296.344546364709)
)

$r = 2
foreach ($row in $forecastRows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Restore the originally-active sheet/tab selection
$wsWeekly.Activate()

Write-Host "PO Forecast sheet populated with" ($r - 2) "data rows"
